# "Set origin to sprite center for flying objects"
# Re-purpose the Tiles sheet data: rename sheet, grow the tile list from
# 44 to 52 rows, and replace the old "+62 / +62 run-length" sequence with
# a (frame, sub-index) pattern used for the flying-object sprite sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Sheet name: "Sheet1" -> "Tiles"
# ---------------------------------------------------------------
$ws.Name = "Tiles"

# ---------------------------------------------------------------
# 2) Grow the sheet from 44 to 52 data rows. Insert 8 rows right after
#    the current last row (44), seeding them from row 44 so the
#    B ("Zarez"/s=1) and C (s=2) column styles carry over exactly like
#    a real Excel "insert copied cells" would.
# ---------------------------------------------------------------
for ($i = 0; $i -lt 8; $i++) {
    $ws.Range("A44:C44").Copy()
    $ws.Range("A45:C45").Insert(-4121)
}

# ---------------------------------------------------------------
# 3) Rows 7-11 get an explicit row height (13.8pt) in the new layout.
# ---------------------------------------------------------------
$ws.Rows("7").RowHeight = 13.8
$ws.Rows("8").RowHeight = 13.8
$ws.Rows("9").RowHeight = 13.8
$ws.Rows("10").RowHeight = 13.8
$ws.Rows("11").RowHeight = 13.8

# ---------------------------------------------------------------
# 4) Column A: first 15 rows become plain literals (no more running
#    "+62" formula chain); from row 16 on it's a simple "+1 on the row
#    above" formula chain.
# ---------------------------------------------------------------
$aValues = @{
    1 = 0
    2 = 1
    3 = 1
    4 = 2
    5 = 2
    6 = 2
    7 = 3
    8 = 3
    9 = 3
    10 = 4
    11 = 4
    12 = 4
    13 = 5
    14 = 5
    15 = 6
}
foreach ($r in $aValues.Keys) {
    $ws.Cells.Item($r, 1).Value = $aValues[$r]
}
for ($r = 16; $r -le 52; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 1).Formula = "=+A$prev+1"
}

# ---------------------------------------------------------------
# 5) Column B: new "sub-index" sequence (1,1,2,1,2,3,1,2,3,1,2,3,1,2,1,
#    1,1,1,...) replacing the constant 62.
# ---------------------------------------------------------------
$bValues = @{
    1 = 1
    2 = 1
    3 = 2
    4 = 1
    5 = 2
    6 = 3
    7 = 1
    8 = 2
    9 = 3
    10 = 1
    11 = 2
    12 = 3
    13 = 1
    14 = 2
    15 = 1
}
foreach ($r in $bValues.Keys) {
    $ws.Cells.Item($r, 2).Value = $bValues[$r]
}
for ($r = 16; $r -le 52; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# ---------------------------------------------------------------
# 6) Column C: constant 4 for every row, no more formula chain.
# ---------------------------------------------------------------
for ($r = 1; $r -le 52; $r++) {
    $ws.Cells.Item($r, 3).Value = 4
}

# ---------------------------------------------------------------
# 7) Selection matches the new data extent.
# ---------------------------------------------------------------
$ws.Range("C2:C52").Select()
